$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "bleu" status label to "noir"
$ws.Cells.Replace("bleu", "noir", 1, 1, $false, $false, $false)

# Correct the wording of the status names (word order fix:
# "... posté" -> "... postés ou publiés")
$ws.Cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés", 1, 1, $false, $false, $false)
$ws.Cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois", 1, 1, $false, $false, $false)
$ws.Cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois", 1, 1, $false, $false, $false)
$ws.Cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés", 1, 1, $false, $false, $false)
